$d = $word.ActiveDocument

$pairs = @(
    @("26×77=2002", "40×72=2880"),
    @("76×43=3268", "48×38=1824"),
    @("18×81=1458", "83×23=1909"),
    @("23×15=345",  "97×42=4074"),
    @("41×63=2583", "92×55=5060"),
    @("78×35=2730", "68×46=3128"),
    @("11×29=319",  "90×66=5940"),
    @("82×76=6232", "24×76=1824"),
    @("47×52=2444", "12×26=312"),
    @("97×57=5529", "83×73=6059"),
    @("83×88=7304", "34×71=2414"),
    @("41×15=615",  "74×45=3330"),
    @("98×43=4214", "85×55=4675"),
    @("98×91=8918", "14×72=1008"),
    @("31×34=1054", "11×87=957"),
    @("60×60=3600", "59×63=3717"),
    @("88×77=6776", "44×21=924"),
    @("67×94=6298", "67×75=5025"),
    @("61×25=1525", "52×58=3016"),
    @("17×94=1598", "16×52=832"),
    @("76×64=4864", "65×72=4680"),
    @("97×65=6305", "96×56=5376"),
    @("51×66=3366", "43×96=4128"),
    @("80×32=2560", "14×13=182"),
    @("26×33=858",  "59×26=1534")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
